$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -12.00939999999999
$ws.Range("B8").Value = 5.891699999999997
$ws.Range("B10").Value = 5.405799999999999
$ws.Range("B12").Value = 5.386599999999998
$ws.Range("C12").Value = -13.8067
$ws.Range("D12").Value = -7.828400000000001
$ws.Range("D13").Value = -8.2799
$ws.Range("C15").Value = -13.4559
$ws.Range("C17").Value = -14.20399999999999
$ws.Range("B18").Value = 6.735099999999996
$ws.Range("D21").Value = -8.458799999999995
$ws.Range("D25").Value = -7.532699999999997
$ws.Range("C26").Value = -12.53610000000001
$ws.Range("C27").Value = -12.57359999999999
$ws.Range("C28").Value = -13.8048
$ws.Range("D32").Value = -7.141700000000002
$ws.Range("D36").Value = -7.444200000000003
$ws.Range("B37").Value = 9.127399999999993
$ws.Range("C37").Value = -13.1836
$ws.Range("D38").Value = -8.0594
$ws.Range("D41").Value = -8.115599999999995
$ws.Range("C47").Value = -12.38839999999999
$ws.Range("D52").Value = -7.879299999999996
$ws.Range("B55").Value = 6.051399999999997
$ws.Range("D59").Value = -8.411899999999992
$ws.Range("C65").Value = -12.4868
$ws.Range("D67").Value = -7.195699999999995
$ws.Range("B68").Value = 4.783499999999996
$ws.Range("C73").Value = -12.0097
$ws.Range("B77").Value = 9.361100000000009
$ws.Range("B78").Value = 10.0454
$ws.Range("B81").Value = 5.599400000000005
$ws.Range("B82").Value = 5.663000000000002
$ws.Range("C84").Value = -13.78129999999999
$ws.Range("D84").Value = -7.920599999999995
$ws.Range("C85").Value = -13.4754
$ws.Range("D88").Value = -7.801299999999998
$ws.Range("D89").Value = -8.259399999999994
$ws.Range("C93").Value = -10.1372
$ws.Range("C95").Value = -13.32109999999999
$ws.Range("D95").Value = -7.611000000000001
$ws.Range("C98").Value = -13.05570000000001
$ws.Range("C99").Value = -12.1364
$ws.Range("C101").Value = -13.427
$ws.Range("D105").Value = -8.191700000000004
